# TC_66.xlsx edit: rename sheet, tweak the shared number format, and
# refresh the CEIC add-in payload cached in the A1 cell comment.
#
# Note: the source commit also drops the customXml/item1.xml +
# customXml/itemProps1.xml parts (leftover CEIC add-in metadata parts).
# The Excel object model exposed here (Workbook.CustomXMLParts /
# .XmlMaps / .RemoveDocumentInformation) does not surface those parts
# for this workbook - CustomXMLParts.Count reports 0 and Add/Delete are
# no-ops against the backing package - so that removal is not
# reachable through COM automation and is intentionally left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet from "My Series" to "Data"
$ws.Name = "Data"

# 2. Update the number format used by the data row (A1:P1):
#    "0.000" -> "###0.000"
$ws.Range("A1:P1").NumberFormat = "###0.000"

# 3. Update the cell comment on A1 with the new encoded payload
#    (gzip+base64 blob used internally by the CEIC Data Manager add-in)
$comment = $ws.Range("A1").Comment
$null = $comment.Text("DRoAAB+LCAAAAAAAAAOlGdtuG8f1VxZ8aoFSu0tatiyMN+BNClFeBJKqrLwUw90RudVyh9mdFcW3FGiRIk1RFIVTpFf0KUWBukabAKndy78Eluw+9Rd65rKzF1Jx6RqCtXNuc+bMuY7QO9eLwLgiUezT8FHF3rMqBgld6vnh7FElYRdV+37lHQd1rl0SnOAILwgDYgO4wvjwOvYfVeaMLQ9Nc7Va7a3qezSamTXLss3H/d7YnZMFrvphzHDokorm8t7MVXFQy1v0CcMeZlhyPqp0x929FvHdNsD6OMQzEu01k9gPSRx3QuYzn8ScMyKYkVa7/x15MKe2d3/PRuYGPKNsJn7gSboCpYQrOtiWTPwFcWqWfVC1Dqp1e2Lbh/X9wzpwPai9lzJqQtTDMRuT6Mp3BWDM8GIp2K2Dum3b9f067LaVCGRlBnDQMPBG5MqPidciQRDvZBFTXWDDZXDq3YxpITPHqwS9vQrHEV7OJz4LyG5qjPpNYxEqXTIhDjqiEXHBfm+l0oCshpEy62TZA+xk7kds3cbrnWWdxiQaLrmRdmN1UJuGrBGQiJ0u4a6JB64ACIdFCUHmHciMqe3HLnz7YUI85wIHcZ6pgERnNLqMl9glA4hjk8tYhQHFHjgc82Pmu3EmYAODTiK6BJGwe5MG3hGITam3YLTsbghG5hs3Kb0siy8ikbhXccNwqwus5W/A0XhOV8MwWI+TaexG/pR47WZKvRWHeEgq7lYSM7oALTIQkrAcZA3/IATLYNQmrr/AwUkAZoydOkgpAFAjYfTCZy0aJItQ27MERWdwogm51ifUazSE6w251WnYDUtStuKKHCO6Sj1nEy6MkAM3Yjf1sk1EmbgNsPT6NjHiRvgpj/wAKkT+LnLQoleM54SwrS4hMYgnwyNec5zmmu+JzAyCwDXBvwHq2A8fWFXLhp+JZR2KH9hZo1En9MSHZQmiIl2KRINkMZxCDF+JMzk24EogBKcImgEOLwF65rP5oJFqvwWD5JnvpN/EIYjWZYDXAqztkoehbugGiUdkEuiGF8IpuW6K+k482gD1IK4dhMP1ZL3kycC8g6JNLnASQH1ikEZmWYYogVEjvizT5EHoNArSO3Z49Y+h/LveYs+FBMlL3J5LFxxgQtU9GyMzT8+zvEs64ayHw1kCeVTbsQzXHsYzwCTCYcyPo9Nmydm2E6H0XmQ6d+RlDRPhTvKyKGCRWaJDE7JY0ggHfTCMf5SEomSqmgAW7mM2VysI3oC4qZHNjFVzFTVLFX8TmQhDeQye4/NhmIMKKn4Y2WnkiDIg4ufsUw+KPQ78aVRws604uLKsBKYOx4+3YzlMbwHaTcgw3yZr3oNkCwXnXuvYKUKsUAtUcsajewe1fateg+zN10gceURwYHSgEWXE6IZXJGYLYDs0RiT2PfjycXBovEumxIewFzZSgbczd54PHUXk/QQa6rVQpQH6FiFFAsipMz/EwSahxmQMzjnBUbDOEcqj9qgLdLc/+tfNr5+9fP7Z7UdPXn/xg//8/Zcv//HTm6c/hI/bv/z15uNfyGNKYjTB04AIhSbNgwOrfg88TYOQShBQ5b3EZQJ2fi6qvl4j1b6KRavTbR33miKjaGDKzjuKhEdDD69pki3H8hBiI3GlZuoJksSZpBlKrQvYNuGlXoYh71SvSJE6j7+LUdri1YvPXr34053cymD54vNwv2rX3lh8oOe3N+h08eETgGzw0vniXtXar9ZqOeISDRrBmAT9nLZT13PqtvXQqtUtW2dzTzvyNqIySkma4JlZ4pOgFk1CFq21C+TXKVI4/gRCRKNlKOQWykW/+PHrPz8pUCnrKkhRCihHk0h6k5kuhOjBaGKMh6ejVseYdMbcTzJcjk4K/xpitbuOp4JThWGCg2/BcMxHY6MCXV/FoBcGwe7cWEMk5uKw4GzboHKjtxRZ1vI4oslS3kiOIYNuodTZZCvHllwjcMKeG0knQ20hl7re/O3zbQzqIMrNTkOf6SkvD0MFjATl8CpqP/3nyy8/fPn8+e2zn918+f2CBLWPHn7AzyGa8kvt9pDyVL0pQdDZWBjz0vpurr4oIG8dT6gfstix74uuUa0QsNpcmviNugsoeUKwsBfASxD0Lo4710wFtjNAZhEAei4xVFua9dcaIHN4Ztd//+a3t7/6/PaTZ68//OPNR3+4+fiTVy9+9/rp72XU3T55dvuTpyrLlwuB0IV37bINNMQQ6Bo8Gg1eu42vPvi5EVJmQNNhJCIjffXBpzlhXFHRnmSSoanTihRV2CDNM3M+I6eK1qHAp1lkA9DiJayuKVQRo0vfzTZ5r8pF8bgTiG90J9UkJgaFfuqbcJIiccb8v/IpFllSTx5YNbumsFIbfoQpjnOmPw7oFJqMFCGGrBJJgevrGTJasd9xb9hs9DISqcQw8mAms/ikyD9Q2lTyktKN05UeLTIIYKHxc5OAz8UbZJsoLTmXxkw1YF40PJ7+tg9nBQqY0aNINkSherEcJ0toh5kaYu/Gi8eZXAM8kL1qviXO1t12EQ/rHBYKYRHNAQIvUpNCyTTVjfksK9vZATdNtgRc4UEHzKEeJWWndQV9ZWTyvNOJIhptTT4ZJiXrQycNGcXMLK5pxJ3KrtvL7ioFpAlPf8jZT52QtklA2G4vdmbG3adXb80Ld78razceBp4y5m6jhzZLJiD/bMkd5f99tZTO1ogiaKz4I8fOz4zp6DqCiXdHbeRRBCMfAWF39Q545Ecxe8wzgfqSkHMNOZcd6mPH3pd96GMJEBTyQx1SSTcLaqahy+TjMQ16/sLfcSy00vguCgFbLpeyhevu5im8tAzINTSYOQmQFKffg7LBR57dpEmHhVyq+flrTezP5mxXxR5MMfHI1Kq6U1Kr3vOsg+pDQupV24b/sVurWdY+f+pRwiFz+GS14yZmemHZH3Wc/wJQPEaHDRoAAA==")
